$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.2095248058516988
$ws.Range("E2").Value = 0.2095248058516988

# Row 3
$ws.Range("D3").Value = 0.9999987493768228
$ws.Range("E3").Value = 0.9999987493768228

# Row 4
$ws.Range("D4").Value = 0.06878562545211682
$ws.Range("E4").Value = 0.06878562545211682

# Row 5
$ws.Range("D5").Value = 0.000000000167084073528315189396
$ws.Range("E5").Value = 0.000000000167084073528315189396

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.1298108064916564
$ws.Range("E6").Value = 0.1298108064916564

# Row 7
$ws.Range("D7").Value = 0.02311207005339746
$ws.Range("E7").Value = 0.9768879299466026

# Row 9
$ws.Range("D9").Value = 0.999641608354239
$ws.Range("E9").Value = 0.0003583916457610137

# Row 10
$ws.Range("D10").Value = 0.06797437125971145
$ws.Range("E10").Value = 0.9320256287402886

# Row 11
$ws.Range("D11").Value = 0.01104048304710489
$ws.Range("E11").Value = 0.9889595169528951
$ws.Range("F11").Value = 2.499987125396729
$ws.Range("G11").Value = 0.6
